# Atualização automática da planilha
# Budget sheet: remove the "Deloitte" line item (row 6) entirely — its
# contents shift the rows below it up by one — and clear out the
# "Realizado" (actual spend) figures for months 1-3 (columns E:G) on every
# remaining budget line except the first two (rows 3 and 4), which keep
# their historical actuals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Budget")
if ($null -eq $ws) { $ws = $wb.ActiveSheet }

# Remove the entire "Consultoria & Implantação / Deloitte / P02" row.
# This shifts rows 7:29 up to 6:28, carrying their formatting with them.
$ws.Rows("6:6").Delete()

# Clear the "Realizado" (E:G, months 1-3) values for the remaining budget
# lines (rows 5-10 after the shift) while keeping their cell formatting.
$ws.Range("E5:G10").ClearContents()

# Reflect the edit location: whole-row selection on the row that used to
# hold the deleted entry (now "Treinamentos / SoftExpert / P01").
$ws.Rows("6:6").Select() | Out-Null
